# Apply updated odds values to Sheet1, matching the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("I2").Value = 19
$ws.Range("U2").Value = 1.92

# Row 4
$ws.Range("V4").Value = 1.41
$ws.Range("AA4").Value = 65
$ws.Range("AB4").Value = 8.6

# Row 5
$ws.Range("F5").Value = 1.64
$ws.Range("G5").Value = 1.66
$ws.Range("H5").Value = 4.7
$ws.Range("I5").Value = 5.2
$ws.Range("J5").Value = 4.9
$ws.Range("R5").Value = 1.9
$ws.Range("S5").Value = 1.9
$ws.Range("U5").Value = 2.72
$ws.Range("W5").Value = 2.5
$ws.Range("X5").Value = 40
$ws.Range("Z5").Value = 55
$ws.Range("AA5").Value = 120
$ws.Range("AC5").Value = 13.5
$ws.Range("AD5").Value = 21
$ws.Range("AE5").Value = 50
$ws.Range("AF5").Value = 16.5
$ws.Range("AG5").Value = 11.5
$ws.Range("AH5").Value = 16.5
$ws.Range("AI5").Value = 44
$ws.Range("AK5").Value = 16
$ws.Range("AL5").Value = 22
$ws.Range("AM5").Value = 55

# Row 6
$ws.Range("T6").Value = 1.72
$ws.Range("V6").Value = 2.04

# Row 7
$ws.Range("L7").Value = 1.01
$ws.Range("M7").Value = 1.01
$ws.Range("N7").Value = 1.01
$ws.Range("O7").Value = 1.01
$ws.Range("P7").Value = 1.08
$ws.Range("R7").Value = 1.08
$ws.Range("S7").Value = 1.01
$ws.Range("T7").Value = 1.01
$ws.Range("U7").Value = 1.01
$ws.Range("V7").Value = 1.01
$ws.Range("W7").Value = 1.01
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 1000
$ws.Range("Z7").Value = 1000
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AE7").Value = 1000
$ws.Range("AF7").Value = 1000
$ws.Range("AG7").Value = 1000
$ws.Range("AH7").Value = 1000
$ws.Range("AI7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("AK7").Value = 1000
$ws.Range("AL7").Value = 1000
$ws.Range("AM7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 1000

# Row 8
$ws.Range("Q8").Value = 2.1
$ws.Range("X8").Value = 13
$ws.Range("Z8").Value = 90
$ws.Range("AA8").Value = 500
$ws.Range("AD8").Value = 38
$ws.Range("AE8").Value = 220
$ws.Range("AG8").Value = 10.5
$ws.Range("AI8").Value = 220
$ws.Range("AJ8").Value = 11.5
$ws.Range("AK8").Value = 17.5
$ws.Range("AL8").Value = 50
$ws.Range("AM8").Value = 300
$ws.Range("AN8").Value = 8.8
$ws.Range("AO8").Value = 390

# Row 9
$ws.Range("H9").Value = 6.2
$ws.Range("I9").Value = 6.4
$ws.Range("P9").Value = 2.14
$ws.Range("Q9").Value = 1.84
$ws.Range("S9").Value = 3.1
$ws.Range("U9").Value = 2.08
$ws.Range("X9").Value = 18
$ws.Range("AC9").Value = 10.5
$ws.Range("AE9").Value = 90
$ws.Range("AF9").Value = 9
$ws.Range("AG9").Value = 9.6
$ws.Range("AH9").Value = 23
$ws.Range("AI9").Value = 95
$ws.Range("AK9").Value = 16.5
$ws.Range("AL9").Value = 36
